# feat: add telegram bot
# Mark rows 6-8 (спорт, семья, проекты) with a "1" flag in column I,
# matching the existing flags already present for rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1
